$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
$rng = $ws.Range("C2:C$lastRow")

foreach ($cell in $rng.Cells) {
    if ($cell.Value2 -eq 45186) {
        $cell.Value2 = 45188
    }
}
